$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("IDName")
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Marks"

$newSheet.Range("A1").Value = "Name"
$newSheet.Range("B1").Value = "Marks"
$newSheet.Range("A2").Value = "Amritha"
$newSheet.Range("B2").Value = 45
$newSheet.Range("A3").Value = "Lakshmi"
$newSheet.Range("B3").Value = 46
$newSheet.Range("A4").Value = "Vyshnavi"
$newSheet.Range("B4").Value = 47

$newSheet.Range("C4").Select()
